# Fruta / hortaliza, semanal
# Insert a new data row at row 33 (pushing existing rows 33-45 down to 34-46)
# and populate it with the latest weekly observation.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows.Item(33).Insert()

$ws.Cells.Item(33, 1).Value = 6
$ws.Cells.Item(33, 2).Value = "Mercado Mayorista Lo Valledor de Santiago"
$ws.Cells.Item(33, 3).Value = "Metropolitana"
$ws.Cells.Item(33, 4).Value = 44452
$ws.Cells.Item(33, 5).Value = 13
$ws.Cells.Item(33, 6).Value = 100114007
$ws.Cells.Item(33, 7).Value = "Jengibre"
$ws.Cells.Item(33, 8).Value = "Sin especificar"
$ws.Cells.Item(33, 9).Value = "Primera"
$ws.Cells.Item(33, 10).Value = 290
$ws.Cells.Item(33, 11).Value = 12000
$ws.Cells.Item(33, 12).Value = 13000
$ws.Cells.Item(33, 13).Value = 12414
$ws.Cells.Item(33, 14).Value = "`$/caja 13 kilos"
$ws.Cells.Item(33, 15).Value = "Perú"
$ws.Cells.Item(33, 16).Value = 955
$ws.Cells.Item(33, 17).Value = 13
$ws.Cells.Item(33, 18).Value = "Hortaliza"
